$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, matching style of existing headers (G1)
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("H1").Value = "Save"

# Add the corresponding numeric value in H2
$ws.Range("H2").Value = 0
